$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 9; this pushes the existing
# rows 9, 10, 11 down to become rows 11, 12, 13 (content unchanged).
$ws.Rows.Item(9).Resize(2).EntireRow.Insert()

# New row 9: fresh weekly price record
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44966
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = "Frutos de pepita"
$ws.Range("I9").Value = 100104005
$ws.Range("J9").Value = "Pera asiática"
$ws.Range("K9").Value = "Hosui"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 250000
$ws.Range("O9").Value = 250000
$ws.Range("P9").Value = 250000
$ws.Range("Q9").Value = "$/bins (400 kilos)"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 625
$ws.Range("T9").Value = 400

# New row 10: fresh weekly price record
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44966
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100104
$ws.Range("H10").Value = "Frutos de pepita"
$ws.Range("I10").Value = 100104005
$ws.Range("J10").Value = "Pera asiática"
$ws.Range("K10").Value = "Hosui"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = "$/caja 18 kilos granel"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 833
$ws.Range("T10").Value = 18

# The previously-existing rows 9, 10, 11 (now shifted to rows 11, 12, 13
# by the insert above) keep their original values unchanged.
